$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 (Iowa) - fill in data that previously failed to parse.
$ws.Range("B42").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B42").Value = (Get-Date -Year 2020 -Month 7 -Day 7 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("C42").Value = 31906
$ws.Range("D42").Value = 725
$ws.Range("E42").Value = 2868
$ws.Range("F42").Value = 36
$ws.Range("G42").Value = 8.99
$ws.Range("H42").Value = 4.97

$ws.Range("I42").Value = $true
$ws.Range("J42").Value = $true

$ws.Range("O42").Value = "Success!"
